$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark that currently sits at the end of the
#    paragraph holding the first screenshot image.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append the new content after the final paragraph
#    (" Este código retorna un objeto creado a parte de un array.") and
#    before the section properties: an empty paragraph, a bold
#    "PARA TENER EN CUENTA" heading paragraph, and the explanatory
#    paragraph that ends with the relocated _GoBack bookmark.
$newContentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PARA TENER EN CUENTA</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Es importante tener en cuenta que la inicialización de los estados </w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>this.state</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = { data: valor}) </w:t></w:r><w:r><w:t xml:space="preserve">solo se puede realizar en el constructor del componentes, a partir de allí si se desea actualizar algún estado solo se puede utilizar la función </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>this.</w:t></w:r><w:r><w:t>setState</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>({ data:nuevoValor</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> }).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)
[void]$insertionRange.InsertXML($newContentXml)
